$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly update swaps the data between the "bandeja 18 kilos" record
# (originally in row 2/3) and the "caja 15 kilos" record (originally in
# row 11/12) for each quality grade (Primera / Segunda), so that the
# newest observation (44536) now sits in rows 11/12 and the older one
# (44424) sits in rows 2/3.
#
# Note: this runtime's Range.Value getter is unreliable, so Value2 is
# used for both reads and writes.

function Swap-RowData($rowA, $rowB) {
    $cols = @("D", "J", "K", "L", "M", "N", "P", "Q")
    foreach ($col in $cols) {
        $rangeA = $ws.Range("$col$rowA")
        $rangeB = $ws.Range("$col$rowB")
        $valA = $rangeA.Value2
        $valB = $rangeB.Value2
        $rangeA.Value2 = $valB
        $rangeB.Value2 = $valA
    }
}

Swap-RowData 2 11
Swap-RowData 3 12
